# Update the "Drive" worksheet of the Drive Summary workbook with the
# latest Qualifying Fuel Spend figures (Operational excluding DQ mapper
# extraction).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drive")
$ws.Activate()

# "Qualifying Fuel Spend" block, "Actual" row (row 6):
#   Km  -> F6
#   Points -> G6
$ws.Range("F6").Value = 330.497
$ws.Range("G6").Value = 192

# "Qualifying Fuel Spend" figure used for the Fuel Cash Back calc (row 13)
$ws.Range("C13").Value = 0

# Leave the cursor where the author left it when saving
$ws.Range("C16").Select()
